$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.0131
$ws.Range("B9").Value = 8.378200000000005
$ws.Range("C11").Value = -13.2373
$ws.Range("B18").Value = 4.858200000000005
$ws.Range("B20").Value = 5.521099999999997
$ws.Range("E21").Value = 13.1169
